$d = $word.ActiveDocument

$searchText       = "Biträdande programansvarig på programmet Digitala upplevelser för lärande. "
$bitradandeBody   = "✅ Biträdande programansvarig på programmet Digitala upplevelser för lärande."
$bitradandeTitle  = "Biträdande programansvarig"
$nominationBody   = "✅ Nominering årets lärare - nominerad till årets lärare 2022 av Blekinge studentkår (BSK). Motiveringen löd `"Strukturerad, hjälpsam och finns alltid där vid snabb återkoppling.`" "
$nominationTitle  = "Nominering årets lärare"

# --- First occurrence: bulleted list item ("ListBullet" style) ---------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found1) {
    throw "Could not find first occurrence of Biträdande programansvarig text"
}

# Trim the trailing space from the existing run's text.
$trim1 = $d.Range($rng1.End - 1, $rng1.End)
$trim1.Text = ""

# Insert a new paragraph right after it, inheriting the ListBullet style,
# and fill it with the nomination bullet text.
$para1 = $rng1.Paragraphs(1)
$para1.Range.InsertParagraphAfter()
$newPara1 = $para1.Next()
$newPara1.Range.InsertBefore($nominationBody)

# --- Second occurrence: plain paragraph followed by a bold heading -----
$rngAfter1 = $d.Range($rng1.End, $d.Content.End)
$found2 = $rngAfter1.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found2) {
    throw "Could not find second occurrence of Biträdande programansvarig text"
}

$para2 = $rngAfter1.Paragraphs(1)

# The bold "Biträdande programansvarig" heading paragraph immediately
# follows the plain-text paragraph just located.
$headingPara = $para2.Next()

if ($headingPara.Range.Text.TrimEnd() -ne $bitradandeTitle) {
    throw "Unexpected paragraph after second occurrence: [$($headingPara.Range.Text)]"
}

# Replace the two paragraphs (plain text + bold heading) with four plain
# paragraphs in one shot so none of the new runs inherit bold formatting;
# bold is re-applied explicitly afterwards only to the two heading runs,
# navigated to directly via the Paragraphs chain (so the similarly-worded
# bullet-list paragraphs earlier in the document are left untouched).
$wholeRange = $d.Range($para2.Range.Start, $headingPara.Range.End)
$wholeRange.Text = $bitradandeBody + "`r" + $bitradandeTitle + "`r" + $nominationBody + "`r" + $nominationTitle

$newBitradandeBodyPara  = $para2
$newBitradandeTitlePara = $newBitradandeBodyPara.Next()
$newNominationBodyPara  = $newBitradandeTitlePara.Next()
$newNominationTitlePara = $newNominationBodyPara.Next()

# Bold only the run text, not the paragraph mark, so no stray <w:pPr>
# run-properties get written (match the heading runs elsewhere in the doc).
$titleRange1 = $d.Range($newBitradandeTitlePara.Range.Start, $newBitradandeTitlePara.Range.End - 1)
$titleRange1.Font.Bold = 1

$titleRange2 = $d.Range($newNominationTitlePara.Range.Start, $newNominationTitlePara.Range.End - 1)
$titleRange2.Font.Bold = 1
